$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.402.80'
$ws.Range("E2").Value = '  -1.59%  '

$ws.Range("D3").Value = '2.042.23'
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.657'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.55'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.60%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.360'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0742'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.50%  '

$ws.Range("E12").Value = '  -4.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.925'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.85%  '

$ws.Range("D15").Value = '2.345.82'
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.11%  '

$ws.Range("D17").Value = '2.036.08'
$ws.Range("E17").Value = '  -1.42%  '

$ws.Range("D18").Value = '36.338.13'
$ws.Range("E18").Value = '  -1.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.62%  '

$ws.Range("D21").Value = '0.0₃0850'
$ws.Range("E21").Value = '  -5.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.38%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -12.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.32%  '

$ws.Range("E30").Value = '  -3.09%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.41%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0586'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0872'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.53%  '

$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0212'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0892'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.53%  '

$ws.Range("D46").Value = '1.381.32'
$ws.Range("E46").Value = '  +5.54%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.54%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.00%  '

$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.97%  '

$ws.Range("D51").Value = '2.229.22'
$ws.Range("E51").Value = '  -0.50%  '
